$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")
$ws.Range("A1").Value = 33915.575708084776
